$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 14.21164676609349
$ws.Cells.Item(2, 3).Value = 6.341266035141048
$ws.Cells.Item(2, 4).Value = 15.24101696135474
$ws.Cells.Item(2, 5).Value = 16.68608037569345
$ws.Cells.Item(2, 7).Value = 66.05293979238395
$ws.Cells.Item(2, 8).Value = 23.1343022577991
$ws.Cells.Item(2, 10).Value = 9.552488925002425
$ws.Cells.Item(2, 11).Value = 13.80292406624594
$ws.Cells.Item(3, 2).Value = 14.100422706785
$ws.Cells.Item(3, 3).Value = 6.288945967409167
$ws.Cells.Item(3, 4).Value = 15.19053247580511
$ws.Cells.Item(3, 5).Value = 16.63663018733889
$ws.Cells.Item(3, 7).Value = 65.32148462783147
$ws.Cells.Item(3, 8).Value = 23.0462838944038
$ws.Cells.Item(3, 10).Value = 9.564298870808141
$ws.Cells.Item(3, 11).Value = 13.75590965787707
$ws.Cells.Item(4, 2).Value = 14.03682153698649
$ws.Cells.Item(4, 3).Value = 6.259787509132004
$ws.Cells.Item(4, 4).Value = 15.1630372322558
$ws.Cells.Item(4, 5).Value = 16.61012458565988
$ws.Cells.Item(4, 7).Value = 64.87748865673913
$ws.Cells.Item(4, 8).Value = 22.99470735299117
$ws.Cells.Item(4, 10).Value = 9.573047574531794
$ws.Cells.Item(4, 11).Value = 13.73157737830806
$ws.Cells.Item(5, 2).Value = 14.01211211869767
$ws.Cells.Item(5, 3).Value = 6.248666377102639
$ws.Cells.Item(5, 4).Value = 15.15272020066381
$ws.Cells.Item(5, 5).Value = 16.60029938963968
$ws.Cells.Item(5, 7).Value = 64.69798885819475
$ws.Cells.Item(5, 8).Value = 22.97431755979775
$ws.Cells.Item(5, 10).Value = 9.57698905300424
$ws.Cells.Item(5, 11).Value = 13.72281079356689
$ws.Cells.Item(6, 2).Value = 14.00808293188171
$ws.Cells.Item(6, 3).Value = 6.246866104183621
$ws.Cells.Item(6, 4).Value = 15.15106084664631
$ws.Cells.Item(6, 5).Value = 16.59872703393698
$ws.Cells.Item(6, 7).Value = 64.66827370666796
$ws.Cells.Item(6, 8).Value = 22.97096999511149
$ws.Cells.Item(6, 10).Value = 9.577666253389296
$ws.Cells.Item(6, 11).Value = 13.72142473003165
$ws.Cells.Item(7, 2).Value = 14.03648336724085
$ws.Cells.Item(7, 3).Value = 6.259634424924284
$ws.Cells.Item(7, 4).Value = 15.16289449149058
$ws.Cells.Item(7, 5).Value = 16.6099881202681
$ws.Cells.Item(7, 7).Value = 64.87506187273789
$ws.Cells.Item(7, 8).Value = 22.99442981656161
$ws.Cells.Item(7, 10).Value = 9.573099207350257
$ws.Cells.Item(7, 11).Value = 13.7314544864765
$ws.Cells.Item(8, 2).Value = 14.17234130001102
$ws.Cells.Item(8, 3).Value = 6.322620169887341
$ws.Cells.Item(8, 4).Value = 15.22288735667469
$ws.Cells.Item(8, 5).Value = 16.66823344719125
$ws.Cells.Item(8, 7).Value = 65.79975904277707
$ws.Cells.Item(8, 8).Value = 23.10344629127432
$ws.Cells.Item(8, 10).Value = 9.556250109541944
$ws.Cells.Item(8, 11).Value = 13.78577801919283
$ws.Cells.Item(9, 2).Value = 14.47458389001394
$ws.Cells.Item(9, 3).Value = 6.468908572508723
$ws.Cells.Item(9, 4).Value = 15.36798021528909
$ws.Cells.Item(9, 5).Value = 16.81273833527309
$ws.Cells.Item(9, 7).Value = 67.64723108459357
$ws.Cells.Item(9, 8).Value = 23.33643943840462
$ws.Cells.Item(9, 10).Value = 9.535097834100807
$ws.Cells.Item(9, 11).Value = 13.92783024543309
$ws.Cells.Item(10, 2).Value = 14.71646376512406
$ws.Cells.Item(10, 3).Value = 6.589133289976676
$ws.Cells.Item(10, 4).Value = 15.49078940116622
$ws.Cells.Item(10, 5).Value = 16.9368779716728
$ws.Cells.Item(10, 7).Value = 69.01657819173062
$ws.Cells.Item(10, 8).Value = 23.51881074308225
$ws.Cells.Item(10, 10).Value = 9.526813702503878
$ws.Cells.Item(10, 11).Value = 14.05313570892995
$ws.Cells.Item(11, 2).Value = 14.83033982968378
$ws.Cells.Item(11, 3).Value = 6.646324970294041
$ws.Cells.Item(11, 4).Value = 15.55004710320138
$ws.Cells.Item(11, 5).Value = 16.99712922286224
$ws.Cells.Item(11, 7).Value = 69.64025075322645
$ws.Cells.Item(11, 8).Value = 23.60408266365114
$ws.Cells.Item(11, 10).Value = 9.524621277885585
$ws.Cells.Item(11, 11).Value = 14.11450419935777
$ws.Cells.Item(12, 2).Value = 14.87397359244073
$ws.Cells.Item(12, 3).Value = 6.668317655840028
$ws.Cells.Item(12, 4).Value = 15.57296161555317
$ws.Cells.Item(12, 5).Value = 17.02047615771623
$ws.Cells.Item(12, 7).Value = 69.87636920209167
$ws.Cells.Item(12, 8).Value = 23.63669359982641
$ws.Cells.Item(12, 10).Value = 9.524017594961892
$ws.Cells.Item(12, 11).Value = 14.13835329965942
$ws.Cells.Item(13, 2).Value = 14.8645542048175
$ws.Cells.Item(13, 3).Value = 6.663566610767331
$ws.Cells.Item(13, 4).Value = 15.56800565526146
$ws.Cells.Item(13, 5).Value = 17.01542456379932
$ws.Cells.Item(13, 7).Value = 69.82552175571007
$ws.Cells.Item(13, 8).Value = 23.62965620057506
$ws.Cells.Item(13, 10).Value = 9.524137535762694
$ws.Cells.Item(13, 11).Value = 14.13319012562071
$ws.Cells.Item(14, 2).Value = 14.83391959683646
$ws.Cells.Item(14, 3).Value = 6.648127739951668
$ws.Cells.Item(14, 4).Value = 15.55192285196766
$ws.Cells.Item(14, 5).Value = 16.99903941947969
$ws.Cells.Item(14, 7).Value = 69.65967833199711
$ws.Cells.Item(14, 8).Value = 23.6067592356074
$ws.Cells.Item(14, 10).Value = 9.524567073136145
$ws.Cells.Item(14, 11).Value = 14.11645415583695
$ws.Cells.Item(15, 2).Value = 14.81522036846761
$ws.Cells.Item(15, 3).Value = 6.638713940564337
$ws.Cells.Item(15, 4).Value = 15.54213312531094
$ws.Cells.Item(15, 5).Value = 16.98907182651745
$ws.Cells.Item(15, 7).Value = 69.55808283253157
$ws.Cells.Item(15, 8).Value = 23.59277553474547
$ws.Cells.Item(15, 10).Value = 9.524859675315552
$ws.Cells.Item(15, 11).Value = 14.10628180926971
$ws.Cells.Item(16, 2).Value = 14.70909552714123
$ws.Cells.Item(16, 3).Value = 6.585443951545342
$ws.Cells.Item(16, 4).Value = 15.48698394000946
$ws.Cells.Item(16, 5).Value = 16.93301542500361
$ws.Cells.Item(16, 7).Value = 68.97582103951602
$ws.Cells.Item(16, 8).Value = 23.51328350612608
$ws.Cells.Item(16, 10).Value = 9.526988690663348
$ws.Cells.Item(16, 11).Value = 14.04921151200572
$ws.Cells.Item(17, 2).Value = 14.64494671172447
$ws.Cells.Item(17, 3).Value = 6.553388092521442
$ws.Cells.Item(17, 4).Value = 15.45401146460023
$ws.Cells.Item(17, 5).Value = 16.89958612900226
$ws.Cells.Item(17, 7).Value = 68.61869950930303
$ws.Cells.Item(17, 8).Value = 23.46510114564461
$ws.Cells.Item(17, 10).Value = 9.528698397313736
$ws.Cells.Item(17, 11).Value = 14.01530687447538
$ws.Cells.Item(18, 2).Value = 14.60841433929092
$ws.Cells.Item(18, 3).Value = 6.535187302737865
$ws.Cells.Item(18, 4).Value = 15.43536658068385
$ws.Cells.Item(18, 5).Value = 16.88071502809677
$ws.Cells.Item(18, 7).Value = 68.41337274527399
$ws.Cells.Item(18, 8).Value = 23.43760613097862
$ws.Cells.Item(18, 10).Value = 9.529830132956453
$ws.Cells.Item(18, 11).Value = 13.99621795149182
$ws.Cells.Item(19, 2).Value = 14.59610892456612
$ws.Cells.Item(19, 3).Value = 6.529066239645311
$ws.Cells.Item(19, 4).Value = 15.42910908029789
$ws.Cells.Item(19, 5).Value = 16.87438719758145
$ws.Cells.Item(19, 7).Value = 68.34387130851275
$ws.Cells.Item(19, 8).Value = 23.42833462215253
$ws.Cells.Item(19, 10).Value = 9.530238801425977
$ws.Cells.Item(19, 11).Value = 13.98982607256329
$ws.Cells.Item(20, 2).Value = 14.65173805275972
$ws.Cells.Item(20, 3).Value = 6.556776149011541
$ws.Cells.Item(20, 4).Value = 15.45748841058724
$ws.Cells.Item(20, 5).Value = 16.90310792043622
$ws.Cells.Item(20, 7).Value = 68.65670858559608
$ws.Cells.Item(20, 8).Value = 23.47020774154234
$ws.Cells.Item(20, 10).Value = 9.52850104215214
$ws.Cells.Item(20, 11).Value = 14.01887355180707
$ws.Cells.Item(21, 2).Value = 14.84290417907446
$ws.Cells.Item(21, 3).Value = 6.652653602852472
$ws.Cells.Item(21, 4).Value = 15.55663397790641
$ws.Cells.Item(21, 5).Value = 17.00383782299314
$ws.Cells.Item(21, 7).Value = 69.70839326983395
$ws.Cells.Item(21, 8).Value = 23.61347603365213
$ws.Cells.Item(21, 10).Value = 9.524434760491379
$ws.Cells.Item(21, 11).Value = 14.12135350262226
$ws.Cells.Item(22, 2).Value = 14.97080323831027
$ws.Cells.Item(22, 3).Value = 6.71725861730508
$ws.Cells.Item(22, 4).Value = 15.62419298485102
$ws.Cells.Item(22, 5).Value = 17.07275989584109
$ws.Cells.Item(22, 7).Value = 70.3953491379275
$ws.Cells.Item(22, 8).Value = 23.70897126811881
$ws.Cells.Item(22, 10).Value = 9.523097564355037
$ws.Cells.Item(22, 11).Value = 14.19187766416808
$ws.Cells.Item(23, 2).Value = 14.90228403805128
$ws.Cells.Item(23, 3).Value = 6.682608017186345
$ws.Cells.Item(23, 4).Value = 15.5878871974626
$ws.Cells.Item(23, 5).Value = 17.03569654184358
$ws.Cells.Item(23, 7).Value = 70.02879525715848
$ws.Cells.Item(23, 8).Value = 23.65783744104265
$ws.Cells.Item(23, 10).Value = 9.523690490064643
$ws.Cells.Item(23, 11).Value = 14.15391913971782
$ws.Cells.Item(24, 2).Value = 14.64866660149077
$ws.Cells.Item(24, 3).Value = 6.5552436949306
$ws.Cells.Item(24, 4).Value = 15.45591551230246
$ws.Cells.Item(24, 5).Value = 16.90151463412931
$ws.Cells.Item(24, 7).Value = 68.63952470989412
$ws.Cells.Item(24, 8).Value = 23.46789840793717
$ws.Cells.Item(24, 10).Value = 9.528589802916589
$ws.Cells.Item(24, 11).Value = 14.01725979939711
$ws.Cells.Item(25, 2).Value = 14.3891848656104
$ws.Cells.Item(25, 3).Value = 6.427007395448777
$ws.Cells.Item(25, 4).Value = 15.32584049440807
$ws.Cells.Item(25, 5).Value = 16.77045025299481
$ws.Cells.Item(25, 7).Value = 67.14469007587725
$ws.Cells.Item(25, 8).Value = 23.27140304451082
$ws.Cells.Item(25, 10).Value = 9.539546133317675
$ws.Cells.Item(25, 11).Value = 13.88566636673001
